$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.359.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "'2.092.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'251.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").Value = "'0.666"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'54.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.90%  "
$ws.Range("D9").Value = "'62.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").Value = "'0.380"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("D11").Value = "'0.0752"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("E12").Value = "  +7.78%  "
$ws.Range("D13").Value = "'15.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.05%  "
$ws.Range("D14").Value = "'2.394.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("E15").Value = "  +3.23%  "
$ws.Range("D16").Value = "'2.093.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "'5.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.93%  "
$ws.Range("D18").Value = "'37.349.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'73.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("D20").Value = "'14.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.27%  "
$ws.Range("D21").Value = "'0.0₃0852"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.08%  "
$ws.Range("D22").Value = "'241.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Value = "'5.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.33%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'172.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.78%  "
$ws.Range("D28").Value = "'20.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("D29").Value = "'2.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.39%  "
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("E31").Value = "  +28.11%  "
$ws.Range("D32").Value = "'22.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.35%  "
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").Value = "'0.0622"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.57%  "
$ws.Range("D35").Value = "'0.0902"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'4.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "'2.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "'1.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("E41").Value = "  +134.10%  "
$ws.Range("E42").Value = "  +13.86%  "
$ws.Range("D43").Value = "'0.0229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.88%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0969"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.12%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'99.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "'1.330.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "'2.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.48%  "
$ws.Range("E50").Value = "  +7.74%  "
$ws.Range("D51").Value = "'6.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.86%  "
